$d = $word.ActiveDocument

# 1. "Allo" -> "Al" (addressee table, first cell)
$d.Content.Find.Execute("Allo", $true, $false, $false, $false, $false, $true, 1, $false, "Al", 2)

# 2. "Sportello unico per le attività produttive" -> "SUAP/SUE" (addressee table, second cell)
$d.Content.Find.Execute("Sportello unico per le attività produttive", $true, $false, $false, $false, $false, $true, 1, $false, "SUAP/SUE", 2) | Out-Null

# 3. " SUAP " -> " SUAP/SUE " (In risposta a nota SUAP <fld> prot.)
$rng3 = $d.Content
$rng3.Find.Execute("In risposta a nota") | Out-Null
$rng3.Collapse(0)
$rng3.Find.Execute("SUAP", $true, $false, $false, $false, $false, $true, 1, $false, "SUAP/SUE", 1) | Out-Null

# 4. "pratica SUAP n°" -> "pratica SUAP/SUE n°"
$d.Content.Find.Execute("pratica SUAP n°", $true, $false, $false, $false, $false, $true, 1, $false, "pratica SUAP/SUE n°", 1) | Out-Null

# 5. Remove the whole "Si premette che il SUAP riceve..." paragraph
$rng5 = $d.Content
$rng5.Find.Execute("Si premette che il")
$para5 = $rng5.Paragraphs(1)
$para5.Range.Delete()
